# Update the cached "datetimeFigureOut" date text (11/8/19 -> 11/9/19)
# wherever it appears: once on the Slide Master, and once on every
# slide layout's Date Placeholder.
$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq "11/8/19") {
                    $shp.TextFrame.TextRange.Text = "11/9/19"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Widen the "Transfer Package" label and rename it to
# "Logic, Command, Transfer Package" on slide 1.
$s = $p.Slides.Item(1)
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Rectangle 28") {
        $shp = $candidate
    }
}
if ($shp -eq $null) {
    $shp = $s.Shapes.Item(3)
}
$shp.Width = [math]::Round(3470181 / 12700, 4)
$shp.TextFrame.TextRange.Text = "Logic, Command, Transfer Package"
